# Append 15 new applicant rows (97-111) to the "Лист1" sheet of qabul.xlsx.
# Columns: A=F.I.SH  B=Passport  C=Shartnoma raqam  D=Viloyat  E=Tuman
#          F=Ta'lim yo'nalishi  G=Telefon raqam  H=Sana
#
# For numeric-looking text (contract numbers in C, phone numbers in G, and
# the one dd-mm date in H111 that Excel could otherwise parse as a date)
# the cell is pre-formatted as Text ("@") before the value is written so it
# is stored as a string instead of being auto-coerced to a number/date, then
# ClearFormats() removes the temporary number format again so the cell is
# left with no explicit style - matching the rest of the sheet's cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97
$ws.Range('A97').Value = 'Ibadullayeva Ozodaxon Nuraddin qizi'
$ws.Range('B97').Value = 'AD4790061'
$ws.Range('C97').NumberFormat = '@'
$ws.Range('C97').Value = '143'
$ws.Range('C97').ClearFormats()
$ws.Range('D97').Value = 'Xorazm viloyati'
$ws.Range('E97').Value = 'Xiva tumani'
$ws.Range('F97').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G97').NumberFormat = '@'
$ws.Range('G97').Value = '998914278764'
$ws.Range('G97').ClearFormats()
$ws.Range('H97').Value = '15-05-2024'

# Row 98
$ws.Range('A98').Value = 'Rustamova Shohnoza Xushmatovna'
$ws.Range('B98').Value = 'AB5924952'
$ws.Range('C98').NumberFormat = '@'
$ws.Range('C98').Value = '144'
$ws.Range('C98').ClearFormats()
$ws.Range('D98').Value = 'Toshkent shahri'
$ws.Range('E98').Value = 'Chilonzor tumani'
$ws.Range('F98').Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi'
$ws.Range('G98').NumberFormat = '@'
$ws.Range('G98').Value = '998999319041'
$ws.Range('G98').ClearFormats()
$ws.Range('H98').Value = '15-05-2024'

# Row 99
$ws.Range('A99').Value = 'Umarova Feruza Isroiljonovna'
$ws.Range('B99').Value = 'AD2736984'
$ws.Range('C99').NumberFormat = '@'
$ws.Range('C99').Value = '145'
$ws.Range('C99').ClearFormats()
$ws.Range('D99').Value = 'Namangan viloyati'
$ws.Range('E99').Value = 'Uchqoʻrgʻon tumani'
$ws.Range('F99').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G99').NumberFormat = '@'
$ws.Range('G99').Value = '998939422615'
$ws.Range('G99').ClearFormats()
$ws.Range('H99').Value = '16-05-2024'

# Row 100
$ws.Range('A100').Value = 'Primqulova Zulfiya Djabbarovna'
$ws.Range('B100').Value = 'AB2914174'
$ws.Range('C100').NumberFormat = '@'
$ws.Range('C100').Value = '146'
$ws.Range('C100').ClearFormats()
$ws.Range('D100').Value = 'Surxondaryo viloyati'
$ws.Range('E100').Value = 'Termiz tumani'
$ws.Range('F100').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G100').NumberFormat = '@'
$ws.Range('G100').Value = '998977940013'
$ws.Range('G100').ClearFormats()
$ws.Range('H100').Value = '16-05-2024'

# Row 101
$ws.Range('A101').Value = 'Xasanova maftuna'
$ws.Range('B101').Value = 'AB1676320'
$ws.Range('C101').NumberFormat = '@'
$ws.Range('C101').Value = '147'
$ws.Range('C101').ClearFormats()
$ws.Range('D101').Value = 'Andijon viloyati'
$ws.Range('E101').Value = 'Andijon tuman'
$ws.Range('F101').Value = 'Maktabgacha ta’lim tashkiloti defektologi/logopedi'
$ws.Range('G101').NumberFormat = '@'
$ws.Range('G101').Value = '998916102810'
$ws.Range('G101').ClearFormats()
$ws.Range('H101').Value = '17-05-2024'

# Row 102
$ws.Range('A102').Value = 'Safarova Dildora Shomurotovna'
$ws.Range('B102').Value = 'AB5494666'
$ws.Range('C102').NumberFormat = '@'
$ws.Range('C102').Value = '148'
$ws.Range('C102').ClearFormats()
$ws.Range('D102').Value = 'Surxondaryo viloyati'
$ws.Range('E102').Value = 'Termiz tumani'
$ws.Range('F102').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G102').NumberFormat = '@'
$ws.Range('G102').Value = '998907478616'
$ws.Range('G102').ClearFormats()
$ws.Range('H102').Value = '17-05-2024'

# Row 103
$ws.Range('A103').Value = 'Xudoyberdiyeva Barchinoy Sobirovna'
$ws.Range('B103').Value = 'AA6486701'
$ws.Range('C103').NumberFormat = '@'
$ws.Range('C103').Value = '149'
$ws.Range('C103').ClearFormats()
$ws.Range('D103').Value = 'Namangan viloyati'
$ws.Range('E103').Value = 'Uchqoʻrgʻon tumani'
$ws.Range('F103').Value = 'Maktabgacha ta’lim tashkiloti musiqa rahbari'
$ws.Range('G103').NumberFormat = '@'
$ws.Range('G103').Value = '998947271989'
$ws.Range('G103').ClearFormats()
$ws.Range('H103').Value = '20-05-2024'

# Row 104
$ws.Range('A104').Value = 'Karimova Umida Hamroboyevna'
$ws.Range('B104').Value = 'AA5316053'
$ws.Range('C104').NumberFormat = '@'
$ws.Range('C104').Value = '150'
$ws.Range('C104').ClearFormats()
$ws.Range('D104').Value = 'Navoiy viloyati'
$ws.Range('E104').Value = 'Nurota tumani'
$ws.Range('F104').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G104').NumberFormat = '@'
$ws.Range('G104').Value = '+998934315343'
$ws.Range('G104').ClearFormats()
$ws.Range('H104').Value = '20-05-2024'

# Row 105
$ws.Range('A105').Value = 'Allaberganova Sharofat Komiljonovna'
$ws.Range('B105').Value = 'AD4771702'
$ws.Range('C105').NumberFormat = '@'
$ws.Range('C105').Value = '151'
$ws.Range('C105').ClearFormats()
$ws.Range('D105').Value = 'Xorazm viloyati'
$ws.Range('E105').Value = 'Urganch tumani'
$ws.Range('F105').Value = 'Maktabgacha ta’lim tashkiloti metodisti'
$ws.Range('G105').NumberFormat = '@'
$ws.Range('G105').Value = '998907192557'
$ws.Range('G105').ClearFormats()
$ws.Range('H105').Value = '20-05-2024'

# Row 106
$ws.Range('A106').Value = 'Yuldasheva Zarifaxon Yashinjan qizi'
$ws.Range('B106').Value = 'AA6363827'
$ws.Range('C106').NumberFormat = '@'
$ws.Range('C106').Value = '152'
$ws.Range('C106').ClearFormats()
$ws.Range('D106').Value = 'Andijon viloyati'
$ws.Range('E106').Value = 'Andijon tuman'
$ws.Range('F106').Value = 'Maktabgacha ta’lim tashkiloti defektologi/logopedi'
$ws.Range('G106').NumberFormat = '@'
$ws.Range('G106').Value = '+79098312494'
$ws.Range('G106').ClearFormats()
$ws.Range('H106').Value = '21-05-2024'

# Row 107
$ws.Range('A107').Value = 'Allayorova Dinora Farhod qizi'
$ws.Range('B107').Value = 'AB3135777'
$ws.Range('C107').NumberFormat = '@'
$ws.Range('C107').Value = '153'
$ws.Range('C107').ClearFormats()
$ws.Range('D107').Value = 'Navoiy viloyati'
$ws.Range('E107').Value = 'Nurota tumani'
$ws.Range('F107').Value = 'Maktabgacha ta’lim tashkiloti musiqa rahbari'
$ws.Range('G107').NumberFormat = '@'
$ws.Range('G107').Value = '+998999573907'
$ws.Range('G107').ClearFormats()
$ws.Range('H107').Value = '22-05-2024'

# Row 108
$ws.Range('A108').Value = 'Badalova Dildora Fayzulla qizi'
$ws.Range('B108').Value = 'AB1866845'
$ws.Range('C108').NumberFormat = '@'
$ws.Range('C108').Value = '154'
$ws.Range('C108').ClearFormats()
$ws.Range('D108').Value = 'Toshkent shahri'
$ws.Range('E108').Value = 'Mirzo Ulugʻbek tumani'
$ws.Range('F108').Value = 'Maktabgacha ta’lim tashkiloti tarbiyachisi'
$ws.Range('G108').NumberFormat = '@'
$ws.Range('G108').Value = '998909996877'
$ws.Range('G108').ClearFormats()
$ws.Range('H108').Value = '23-05-2024'

# Row 109
$ws.Range('A109').Value = 'Abdullayeva Minajaat Mahkamovna'
$ws.Range('B109').Value = 'AB6360027'
$ws.Range('C109').NumberFormat = '@'
$ws.Range('C109').Value = '155'
$ws.Range('C109').ClearFormats()
$ws.Range('D109').Value = 'Namangan viloyati'
$ws.Range('E109').Value = 'Uchqoʻrgʻon tumani'
$ws.Range('F109').Value = 'Maktabgacha ta’lim tashkiloti direktori'
$ws.Range('G109').NumberFormat = '@'
$ws.Range('G109').Value = '998937092759'
$ws.Range('G109').ClearFormats()
$ws.Range('H109').Value = '27-05-2024'

# Row 110
$ws.Range('A110').Value = 'Sulaymonova Feruza Baxshulloevna'
$ws.Range('B110').Value = 'AA6765464'
$ws.Range('C110').NumberFormat = '@'
$ws.Range('C110').Value = '156'
$ws.Range('C110').ClearFormats()
$ws.Range('D110').Value = 'Navoiy viloyati'
$ws.Range('E110').Value = 'Navoiy shahri'
$ws.Range('F110').Value = 'Maktabgacha ta’lim tashkiloti defektologi/logopedi'
$ws.Range('G110').NumberFormat = '@'
$ws.Range('G110').Value = '998913308498'
$ws.Range('G110').ClearFormats()
$ws.Range('H110').Value = '29-05-2024'

# Row 111
$ws.Range('A111').Value = 'Toxtasinova Munojat Muxammadjon qizi'
$ws.Range('B111').Value = 'AC0335486'
$ws.Range('C111').NumberFormat = '@'
$ws.Range('C111').Value = '157'
$ws.Range('C111').ClearFormats()
$ws.Range('D111').Value = 'Namangan viloyati'
$ws.Range('E111').Value = 'Uchqoʻrgʻon tumani'
$ws.Range('F111').Value = 'Maktabgacha ta’lim tashkiloti musiqa rahbari'
$ws.Range('G111').NumberFormat = '@'
$ws.Range('G111').Value = '998933338395'
$ws.Range('G111').ClearFormats()
$ws.Range('H111').NumberFormat = '@'
$ws.Range('H111').Value = '05-06-2024'
$ws.Range('H111').ClearFormats()

